# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the K column (G) values for rows 2-30 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2,5,6,5,5,7,4,4,3,4,2,4,2,5,3,7,1,3,7,1,0,8,5,4,4,4,5,3,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
